$d = $word.ActiveDocument

# --- 1. Split the "Création des éléments de base : " run into three runs
#        ("Créa" / "ti" / "on des éléments de base : ") while keeping the
#        paragraph's own formatting (Titre2 style + both justification)
#        untouched. We locate the paragraph (it is the 8th paragraph of the
#        document) and replace its content via a raw OOXML fragment
#        (InsertXML) so the three <w:r> runs come out completely clean (no
#        stray rPr/rsid noise). NB: string equality against the document's
#        own text is unreliable here (French typography uses non-breaking
#        spaces etc.), so we address paragraphs by position instead.

$target = $d.Paragraphs.Item(8)

$splitXml = @'
<?xml version="1.0" standalone="yes"?>
<pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage"><pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml"><pkg:xmlData><w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:body><w:p><w:r><w:t>Créa</w:t></w:r><w:r><w:t>ti</w:t></w:r><w:r><w:t xml:space="preserve">on des éléments de base : </w:t></w:r></w:p></w:body></w:document></pkg:xmlData></pkg:part></pkg:package>
'@

$start = $target.Range.Start
$end = $target.Range.End
$body = $d.Range($start, $end - 1)
$body.Text = ""
$insertionPoint = $d.Range($start, $start)
$insertionPoint.InsertXML($splitXml)

# --- 2. Append the new "Création de la page d'accueil" section (a Titre2
#        heading paragraph, a justified body paragraph made of several
#        runs, and two trailing empty paragraphs) right after the
#        paragraph that ends the previous section (the one ending in "»"),
#        which is the last (9th) paragraph of the original document.

$lastPara = $d.Paragraphs.Item(9)

$newSectionXml = @'
<?xml version="1.0" standalone="yes"?>
<pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage"><pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml"><pkg:xmlData><w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:body>
<w:p><w:pPr><w:pStyle w:val="Titre2"/></w:pPr><w:r><w:t xml:space="preserve">Création de la page d’accueil : </w:t></w:r></w:p>
<w:p><w:pPr><w:jc w:val="both"/></w:pPr>
<w:r><w:t xml:space="preserve">Je m’atèle ensuite à la création de la page d’accueil </w:t></w:r>
<w:r><w:t>qui conformément à la consigne doit contenir une présentation de l’API ainsi que de ces différents paramètres pour les recherches.</w:t></w:r>
<w:r><w:t xml:space="preserve"> Cette partie du travail est encore une fois du front car il s’agit </w:t></w:r>
<w:r><w:t xml:space="preserve">encore une fois </w:t></w:r>
<w:r><w:t xml:space="preserve">d’ajouter tout simplement le contenu </w:t></w:r>
<w:r><w:t>de la page d’accueil</w:t></w:r>
<w:r><w:t xml:space="preserve"> ce qui permet au passage </w:t></w:r>
<w:r><w:t>d’apprendre les détails de chacun des arguments de l’API.</w:t></w:r>
<w:r><w:t xml:space="preserve"> Cette page à été finis vers les 15h00 </w:t></w:r>
</w:p>
<w:p/>
<w:p/>
</w:body></w:document></pkg:xmlData></pkg:part></pkg:package>
'@

$tailPoint = $d.Range($lastPara.Range.End, $lastPara.Range.End)
$tailPoint.InsertXML($newSectionXml)

Write-Output "done"
